$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 376.05884
$ws.Range("I6").Value = 296.7
$ws.Range("J6").Value = 489.42856
$ws.Range("K6").Value = 890.0999999999999
$ws.Range("L6").Value = 1468.28568
$ws.Range("M6").Value = -778.0999999999999
$ws.Range("N6").Value = -1692.28568
$ws.Range("H52").Value = 1728.5714
$ws.Range("I52").Value = 1750
$ws.Range("K52").Value = 5250
$ws.Range("M52").Value = -5090
$ws.Range("H86").Value = 2958.8096
$ws.Range("I86").Value = 3017.5334
$ws.Range("J86").Value = 2926.1853
$ws.Range("K86").Value = 3017.5334
$ws.Range("L86").Value = 2926.1853
$ws.Range("M86").Value = -1894.5334
$ws.Range("N86").Value = -5172.1853
$ws.Range("H89").Value = 2958.8096
$ws.Range("I89").Value = 3017.5334
$ws.Range("J89").Value = 2926.1853
$ws.Range("K89").Value = 15087.667
$ws.Range("L89").Value = 14630.9265
$ws.Range("M89").Value = -9471.666999999999
$ws.Range("N89").Value = -25862.9265
$ws.Range("H111").Value = 385.3125
$ws.Range("I111").Value = 452.5
$ws.Range("J111").Value = 273.33334
$ws.Range("K111").Value = 1357.5
$ws.Range("L111").Value = 820.0000200000001
$ws.Range("M111").Value = 1709.5
$ws.Range("N111").Value = -6954.00002
$ws.Range("H116").Value = 11367254
$ws.Range("I116").Value = 3232.9333
$ws.Range("J116").Value = 35718730
$ws.Range("K116").Value = 3232.9333
$ws.Range("L116").Value = 35718730
$ws.Range("M116").Value = 209.0666999999999
$ws.Range("N116").Value = -35725614
$ws.Range("H132").Value = 2842135.2
$ws.Range("I132").Value = 2842135.2
$ws.Range("K132").Value = 8526405.600000001
$ws.Range("M132").Value = -8523875.600000001
$ws.Range("H137").Value = 4386.5264
$ws.Range("I137").Value = 4249.3
$ws.Range("J137").Value = 4539
$ws.Range("K137").Value = 12747.9
$ws.Range("L137").Value = 13617
$ws.Range("M137").Value = -10197.9
$ws.Range("N137").Value = -18717

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4370.55
$ws.Range("I32").Value = 3745.3975
$ws.Range("J32").Value = 6587
$ws.Range("K32").Value = 3745.3975
$ws.Range("L32").Value = 6587
$ws.Range("M32").Value = -3458.3975
$ws.Range("N32").Value = -7161
$ws.Range("H45").Value = 1624
$ws.Range("I45").Value = 1604.4445
$ws.Range("J45").Value = 1800
$ws.Range("K45").Value = 1604.4445
$ws.Range("L45").Value = 1800
$ws.Range("M45").Value = -1227.4445
$ws.Range("N45").Value = -2554
$ws.Range("H122").Value = 2501
$ws.Range("I122").Value = 2000.45
$ws.Range("J122").Value = 5838
$ws.Range("K122").Value = 6001.35
$ws.Range("L122").Value = 17514
$ws.Range("M122").Value = -3551.35
$ws.Range("N122").Value = -22414
$ws.Range("H132").Value = 2556.2
$ws.Range("I132").Value = 2233.2954
$ws.Range("J132").Value = 3847.818
$ws.Range("K132").Value = 6699.8862
$ws.Range("L132").Value = 11543.454
$ws.Range("M132").Value = -4169.8862
$ws.Range("N132").Value = -16603.454
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1519.5758
$ws.Range("I134").Value = 1306.1111
$ws.Range("J134").Value = 1775.7333
$ws.Range("K134").Value = 3918.3333
$ws.Range("L134").Value = 5327.199900000001
$ws.Range("M134").Value = -1383.3333
$ws.Range("N134").Value = -10397.1999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2092.3125
$ws.Range("I99").Value = 1487
$ws.Range("J99").Value = 2367.4546
$ws.Range("K99").Value = 1487
$ws.Range("L99").Value = 2367.4546
$ws.Range("M99").Value = 11
$ws.Range("N99").Value = -5363.4546
$ws.Range("H122").Value = 5750.3
$ws.Range("I122").Value = 6404.8
$ws.Range("J122").Value = 5095.8
$ws.Range("K122").Value = 19214.4
$ws.Range("L122").Value = 15287.4
$ws.Range("M122").Value = -16764.4
$ws.Range("N122").Value = -20187.4
$ws.Range("H126").Value = 2092.3125
$ws.Range("I126").Value = 1487
$ws.Range("J126").Value = 2367.4546
$ws.Range("K126").Value = 4461
$ws.Range("L126").Value = 7102.3638
$ws.Range("M126").Value = -1991
$ws.Range("N126").Value = -12042.3638
$ws.Range("H132").Value = 1431.5883
$ws.Range("I132").Value = 968.13794
$ws.Range("J132").Value = 4119.6
$ws.Range("K132").Value = 2904.41382
$ws.Range("L132").Value = 12358.8
$ws.Range("M132").Value = -374.4138199999998
$ws.Range("N132").Value = -17418.8
$ws.Range("H134").Value = 1708.8485
$ws.Range("I134").Value = 1756.8846
$ws.Range("J134").Value = 1530.4286
$ws.Range("K134").Value = 5270.6538
$ws.Range("L134").Value = 4591.2858
$ws.Range("M134").Value = -2735.6538
$ws.Range("N134").Value = -9661.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1711.6666
$ws.Range("I3").Value = 1715.7142
$ws.Range("J3").Value = 1709.091
$ws.Range("K3").Value = 5147.142599999999
$ws.Range("L3").Value = 5127.272999999999
$ws.Range("M3").Value = -5035.142599999999
$ws.Range("N3").Value = -5351.272999999999
$ws.Range("H4").Value = 55624.223
$ws.Range("I4").Value = 61.142857
$ws.Range("J4").Value = 250095
$ws.Range("K4").Value = 183.428571
$ws.Range("L4").Value = 750285
$ws.Range("M4").Value = -71.42857100000001
$ws.Range("N4").Value = -750509
$ws.Range("H12").Value = 52.296295
$ws.Range("I12").Value = 93.454544
$ws.Range("J12").Value = 24
$ws.Range("K12").Value = 280.363632
$ws.Range("L12").Value = 72
$ws.Range("M12").Value = -107.363632
$ws.Range("N12").Value = -418
$ws.Range("H68").Value = 11665.667
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 12998.875
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 38996.625
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -40618.625
$ws.Range("H71").Value = 11665.667
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 12998.875
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 116989.875
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -125101.875
$ws.Range("H74").Value = 9015
$ws.Range("I74").Value = 6740
$ws.Range("J74").Value = 10278.889
$ws.Range("K74").Value = 20220
$ws.Range("L74").Value = 30836.667
$ws.Range("M74").Value = -19159
$ws.Range("N74").Value = -32958.667
$ws.Range("H77").Value = 9015
$ws.Range("I77").Value = 6740
$ws.Range("J77").Value = 10278.889
$ws.Range("K77").Value = 60660
$ws.Range("L77").Value = 92510.00099999999
$ws.Range("M77").Value = -55356
$ws.Range("N77").Value = -103118.001
$ws.Range("H133").Value = 2322.4644
$ws.Range("I133").Value = 2153.16
$ws.Range("J133").Value = 3733.3333
$ws.Range("K133").Value = 6459.48
$ws.Range("L133").Value = 11199.9999
$ws.Range("M133").Value = -1399.48
$ws.Range("N133").Value = -21319.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3801.3333
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws.Range("H99").Value = 1732.1666
$ws.Range("I99").Value = 1732.1666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1732.1666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 513.8334
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2828.0527
$ws.Range("I40").Value = 2591.9167
$ws.Range("K40").Value = 2591.9167
$ws.Range("M40").Value = -2455.9167
$ws.Range("H61").Value = 2473.6667
$ws.Range("I61").Value = 2050.4167
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 2050.4167
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -1848.4167
$ws.Range("N61").Value = -4570.6665
$ws.Range("H113").Value = 2473.6667
$ws.Range("I113").Value = 2050.4167
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 2050.4167
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 119.5832999999998
$ws.Range("N113").Value = -8506.666499999999
$ws.Range("H122").Value = 4215.3335
$ws.Range("I122").Value = 5012.3335
$ws.Range("J122").Value = 3418.3333
$ws.Range("K122").Value = 15037.0005
$ws.Range("L122").Value = 10254.9999
$ws.Range("M122").Value = -12587.0005
$ws.Range("N122").Value = -15154.9999
$ws.Range("H132").Value = 3887.178
$ws.Range("I132").Value = 2527.24
$ws.Range("K132").Value = 7581.719999999999
$ws.Range("M132").Value = -5051.719999999999
$ws.Range("H136").Value = 2903.4902
$ws.Range("I136").Value = 2527.3044
$ws.Range("J136").Value = 3212.5
$ws.Range("K136").Value = 7581.9132
$ws.Range("L136").Value = 9637.5
$ws.Range("M136").Value = -5031.9132
$ws.Range("N136").Value = -14737.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 369.8
$ws.Range("I113").Value = 237.25
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 711.75
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = 1458.25
$ws.Range("N113").Value = -7040
$ws.Range("H122").Value = 668438.9399999999
$ws.Range("I122").Value = 2000820.8
$ws.Range("J122").Value = 2248
$ws.Range("K122").Value = 6002462.4
$ws.Range("L122").Value = 6744
$ws.Range("M122").Value = -6000012.4
$ws.Range("N122").Value = -11644
$ws.Range("H132").Value = 908.5439
$ws.Range("I132").Value = 640.86664
$ws.Range("J132").Value = 1912.3334
$ws.Range("K132").Value = 1922.59992
$ws.Range("L132").Value = 5737.0002
$ws.Range("M132").Value = 607.4000800000001
$ws.Range("N132").Value = -10797.0002
